# Bugfix: Remove quotes around YES in elimination
# Close #118
#
# This adds a new "Data" worksheet (with its own table) containing the
# population-by-place-of-birth/time/gender breakdown, and restyles the
# three existing tables from TableStyleLight9 to TableStyleMedium2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0. The engine auto-names a freshly created table "Table<N+1>" (where N
#    is the current number of tables in the workbook). With 3 existing
#    tables (Table3/Table4/Table5) that means a brand-new table would be
#    auto-named "Table4", colliding with the existing "Table4" ListObject
#    on the Variables sheet. Temporarily rename that one out of the way
#    so the new table can be created and renamed safely, then restore it.
# ---------------------------------------------------------------------
$variablesWs = $wb.Worksheets.Item("Variables")
$variablesLo = $variablesWs.ListObjects.Item(1)
$variablesLo.Name = "TempTable4Holder"

# ---------------------------------------------------------------------
# 1. Add the new "Data" worksheet as the last sheet in the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Data"

# View settings matching the other data sheets
$ws.Application.ActiveWindow.Zoom = 100
$ws.Application.ActiveWindow.DisplayGridlines = $false

# Column widths
$ws.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws.Columns.Item(2).ColumnWidth = 3.8333333333333335
$ws.Columns.Item(3).ColumnWidth = 5.833333333333333
$ws.Columns.Item(4).ColumnWidth = 7.833333333333333

# Page setup
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "place of birth"
$ws.Range("B1").Value = "time"
$ws.Range("C1").Value = "gender"
$ws.Range("D1").Value = "figures_"

# ---------------------------------------------------------------------
# 3. Data rows (place of birth, time, gender, figures_)
#    Values are written as text so that they round-trip through the
#    shared string table the same way the source data used.
# ---------------------------------------------------------------------
$dataRows = @(
    @("T","2018","T","55877"),
    @("T","2018","M","29489"),
    @("T","2018","K","26388"),
    @("T","2019","T","55992"),
    @("T","2019","M","29553"),
    @("T","2019","K","26439"),
    @("T","2020","T","56081"),
    @("T","2020","M","29551"),
    @("T","2020","K","26530"),
    @("T","2021","T","56421"),
    @("T","2021","M","29749"),
    @("T","2021","K","26672"),
    @("T","2022","T","56562"),
    @("T","2022","M","29803"),
    @("T","2022","K","26759"),
    @("N","2018","T","50171"),
    @("N","2018","M","25779"),
    @("N","2018","K","24392"),
    @("N","2019","T","50251"),
    @("N","2019","M","25817"),
    @("N","2019","K","24434"),
    @("N","2020","T","50190"),
    @("N","2020","M","25738"),
    @("N","2020","K","24452"),
    @("N","2021","T","50365"),
    @("N","2021","M","25824"),
    @("N","2021","K","24541"),
    @("N","2022","T","50388"),
    @("N","2022","M","25810"),
    @("N","2022","K","24578"),
    @("S","2018","T","5706"),
    @("S","2018","M","3710"),
    @("S","2018","K","1996"),
    @("S","2019","T","5741"),
    @("S","2019","M","3736"),
    @("S","2019","K","2005"),
    @("S","2020","T","5891"),
    @("S","2020","M","3813"),
    @("S","2020","K","2078"),
    @("S","2021","T","6056"),
    @("S","2021","M","3925"),
    @("S","2021","K","2131"),
    @("S","2022","T","6174"),
    @("S","2022","M","3993"),
    @("S","2022","K","2181")
)

$r = 2
foreach ($row in $dataRows) {
    $rng = $ws.Range("A$r" + ":D$r")
    $rng.NumberFormat = "@"
    $arr = New-Object 'object[,]' 1,4
    for ($i = 0; $i -lt 4; $i++) { $arr[0,$i] = $row[$i] }
    $rng.Value = $arr
    $rng.ClearFormats()
    $r++
}

# ---------------------------------------------------------------------
# 4. Turn the range into a table ("Table6") styled like the others
# ---------------------------------------------------------------------
$lastRow = $r - 1
$tableRange = $ws.Range("A1:D$lastRow")
$lo = $ws.ListObjects.Add(1, $tableRange)
$lo.Name = "Table6"
$lo.TableStyle = "TableStyleMedium2"

# Restore the Variables sheet table's original name
$variablesLo.Name = "Table4"

# ---------------------------------------------------------------------
# 5. Re-style the existing tables (Table3/Table4/Table5) on the first
#    three sheets: TableStyleLight9 -> TableStyleMedium2
# ---------------------------------------------------------------------
foreach ($sheetName in @("Table", "Variables", "Codelists")) {
    $existingWs = $wb.Worksheets.Item($sheetName)
    foreach ($existingLo in $existingWs.ListObjects) {
        $existingLo.TableStyle = "TableStyleMedium2"
    }
}

# ---------------------------------------------------------------------
# 6. Restore "Table" as the active sheet, as it was before the edit
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Table").Activate()
